# Project 1 - code till Sep 25th
# Applies the MasterTestData.xlsx edit:
#  - rename "Sheet1" -> "smoke"
#  - "regression" sheet: drop D1/D2, add rows 7/8 (TC103/Location, TC103/Sydney),
#    widen column B, move the active selection to B7
#  - "smoke" sheet: collapse its old multi-cell selection back to A1

$wb = $excel.ActiveWorkbook

$regression = $wb.Worksheets.Item("regression")
$smoke      = $wb.Worksheets.Item("Sheet1")

# --- "smoke" sheet: clear out its stale selection (was A1:D2) ------------
$smoke.Range("A1").Select()

# --- rename Sheet1 -> smoke ------------------------------------------------
$smoke.Name = "smoke"

# --- "regression" sheet edits ----------------------------------------------
# Drop the old "Expected Title" column entries on rows 1-2
$regression.Range("D1").ClearContents()
$regression.Range("D2").ClearContents()

# Add the two new TC103 rows
$regression.Range("A7").Value = "TC103"
$regression.Range("B7").Value = "Location"
$regression.Range("A8").Value = "TC103"
$regression.Range("B8").Value = "Sydney"

# Widen column B to fit the new/longer content
$regression.Columns.Item(2).ColumnWidth = 56.3

# Put the final selection on the regression sheet (which stays the active tab)
$regression.Range("B7").Select()
